$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Computation time" table (rows 1-4) ---
$ws.Range("B2").Value = 39.421962000000001
$ws.Range("C2").Value = 39.756878
$ws.Range("D2").Value = 40.743558
$ws.Range("E2").Value = 40.180394

$ws.Range("B3").Value = 38.617897999999997
$ws.Range("C3").Value = 38.130916999999997
$ws.Range("D3").Value = 38.902315000000002
$ws.Range("E3").Value = 39.032794000000003

$ws.Range("B4").Value = 20.035744999999999
$ws.Range("C4").Value = 19.886908999999999
$ws.Range("D4").Value = 20.764433
$ws.Range("E4").Value = 20.221336000000001

# --- "Total time" table (rows 9-12) ---
$ws.Range("B10").Value = 0.13971600000000001
$ws.Range("C10").Value = 0.139907
$ws.Range("D10").Value = 0.137294
$ws.Range("E10").Value = 0.13885500000000001

$ws.Range("B11").Value = 0.15398999999999999
$ws.Range("C11").Value = 0.149233
$ws.Range("D11").Value = 0.14890100000000001
$ws.Range("E11").Value = 0.15029400000000001

$ws.Range("B12").Value = 0.47295999999999999
$ws.Range("C12").Value = 0.47224100000000002
$ws.Range("D12").Value = 0.46598099999999998
$ws.Range("E12").Value = 0.46965499999999999

# --- "Overheads time" table (rows 17-20) ---
$ws.Range("B18").Value = 39.282246000000001
$ws.Range("C18").Value = 39.616970999999999
$ws.Range("D18").Value = 40.906264
$ws.Range("E18").Value = 40.041539

$ws.Range("B19").Value = 38.463908000000004
$ws.Range("C19").Value = 37.981684000000001
$ws.Range("D19").Value = 38.753413999999999
$ws.Range("E19").Value = 38.8825

$ws.Range("B20").Value = 19.562784000000001
$ws.Range("C20").Value = 19.414667999999999
$ws.Range("D20").Value = 20.298452000000001
$ws.Range("E20").Value = 19.751681000000001

# --- Update the view: top-left cell and selection ---
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("B18:F20").Select()
